# Update the "practiceQ" worksheet:
#  - new SubmitResult column (G)
#  - replace the buggy findMaxConsecutiveOnes snippet with a fixed version
#  - replace the duplicated snippet in column D with a new findNumbers snippet
#  - fold the old standalone SubmitResult row into row 4 (new column G)
#  - rename "RunResult" label to "Result"
#  - refresh run-result values for the two updated snippets
#  - tweak column widths / row height / alignment to match

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("practiceQ")

$findMaxConsecutiveOnes = "def findMaxConsecutiveOnes(nums):`n    max_count = 0`n    current_count = 0`n    for num in nums:`n        if num == 1:`n            current_count += 1`n            max_count = max(max_count, current_count)`n        else:`n            current_count = 0  # reset when 0 is found`n    return max_count`n# Example usage:`nprint(findMaxConsecutiveOnes([1,1,0,1,1,1]))  # Output: 3`nprint(findMaxConsecutiveOnes([1,0,1,1,0,1]))  # Output: 2"

$findNumbers = "def findNumbers(nums):`n    count = 0`n    for num in nums:`n        if len(str(num)) % 2 == 0:`n            count += 1`n    return count`n# Example usage:`nprint(findNumbers([12, 345, 2, 6, 7896]))   # Output: 2`nprint(findNumbers([555, 901, 482, 1771]))  # Output: 1"

# --- Row 1 header: new SubmitResult column ---
$ws.Range("G1").Value = "SubmitResult"

# --- Row 2: the code samples (and their formatting) ---
$ws.Range("C2").Value = $findMaxConsecutiveOnes
$ws.Range("D2").Value = $findNumbers

$ws.Range("B2:E2").WrapText = $true
$ws.Range("B2:E2").VerticalAlignment = -4160

$ws.Range("F2").WrapText = $false
$ws.Range("F2").VerticalAlignment = -4160

# --- Row 4: rename label, refresh results, pull SubmitResult in as col G ---
$ws.Range("A4").Value = "Result"

$ws.Range("B4").WrapText = $true
$ws.Range("B4").VerticalAlignment = -4160

$ws.Range("C4").Value = "3`n2"
$ws.Range("D4").Value = "2`n1"

$ws.Range("F4").Value = "hello"
$ws.Range("F4").WrapText = $true
$ws.Range("F4").VerticalAlignment = -4160

$ws.Range("G4").Value = "Submission Successful"
$ws.Range("G4").WrapText = $true
$ws.Range("G4").VerticalAlignment = -4160

# Old standalone SubmitResult row is no longer needed - folded into row 4/col G
$ws.Rows.Item(5).Delete()

$ws.Rows.Item(4).RowHeight = 29

# Leftover formatted-but-empty cell below the table
$ws.Range("D6").VerticalAlignment = -4160

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 52.5
$ws.Columns.Item(7).ColumnWidth = 21.67

# --- Selection / view ---
$ws.Range("G4").Select()
